$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.22
$ws.Range("G2").Value = 1.28
$ws.Range("H2").Value = 14
$ws.Range("I2").Value = 26
$ws.Range("J2").Value = 5.6
$ws.Range("K2").Value = 7.6
$ws.Range("L2").Value = 1.27
$ws.Range("N2").Value = 4.3
$ws.Range("P2").Value = 2.14
$ws.Range("Q2").Value = 1.71
$ws.Range("R2").Value = 1.45
$ws.Range("S2").Value = 2.8
$ws.Range("T2").Value = 2.48
$ws.Range("V2").Value = 1.04
$ws.Range("AB2").Value = 9.199999999999999
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 85
$ws.Range("AG2").Value = 14.5
$ws.Range("AH2").Value = 60
$ws.Range("AJ2").Value = 8.800000000000001
$ws.Range("AK2").Value = 19
$ws.Range("AL2").Value = 70
$ws.Range("AN2").Value = 5.1
$ws.Range("F3").Value = 3.05
$ws.Range("G3").Value = 3.4
$ws.Range("H3").Value = 2.46
$ws.Range("I3").Value = 2.64
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 3.6
$ws.Range("L3").Value = 1.45
$ws.Range("N3").Value = 3.5
$ws.Range("O3").Value = 1.35
$ws.Range("P3").Value = 1.84
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.31
$ws.Range("T3").Value = 1.76
$ws.Range("U3").Value = 2.12
$ws.Range("V3").Value = 1.62
$ws.Range("X3").Value = 13.5
$ws.Range("AB3").Value = 12.5
$ws.Range("AC3").Value = 7.8
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 75
$ws.Range("AF3").Value = 23
$ws.Range("AG3").Value = 14
$ws.Range("AH3").Value = 18
$ws.Range("AI3").Value = 120
$ws.Range("AL3").Value = 250
$ws.Range("AN3").Value = 75
$ws.Range("F4").Value = 12.5
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 1.27
$ws.Range("I4").Value = 1.31
$ws.Range("N4").Value = 4.5
$ws.Range("P4").Value = 2.32
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 2.54
$ws.Range("T4").Value = 2.12
$ws.Range("U4").Value = 1.69
$ws.Range("V4").Value = 4.1
$ws.Range("W4").Value = 1.06
$ws.Range("X4").Value = 22
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 10
$ws.Range("AA4").Value = 9.6
$ws.Range("AC4").Value = 14.5
$ws.Range("AD4").Value = 11.5
$ws.Range("AE4").Value = 15.5
$ws.Range("AG4").Value = 55
$ws.Range("AH4").Value = 42
$ws.Range("AI4").Value = 46
$ws.Range("AO4").Value = 5.1
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 17.5
$ws.Range("H5").Value = 1.26
$ws.Range("I5").Value = 1.34
$ws.Range("J5").Value = 5.3
$ws.Range("K5").Value = 7.4
$ws.Range("L5").Value = 1.26
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 5.2
$ws.Range("P5").Value = 2.52
$ws.Range("Q5").Value = 1.48
$ws.Range("S5").Value = 2.26
$ws.Range("T5").Value = 1.91
$ws.Range("V5").Value = 3.9
$ws.Range("W5").Value = 1.06
$ws.Range("AO5").Value = 15
$ws.Range("F6").Value = 2.46
$ws.Range("G6").Value = 2.66
$ws.Range("H6").Value = 2.76
$ws.Range("I6").Value = 2.96
$ws.Range("J6").Value = 3.7
$ws.Range("L6").Value = 1.38
$ws.Range("V6").Value = 1.5
$ws.Range("W6").Value = 1.6
$ws.Range("Y6").Value = 25
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 1000
$ws.Range("AC6").Value = 8.800000000000001
$ws.Range("AD6").Value = 26
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 80
$ws.Range("AG6").Value = 26
$ws.Range("AJ6").Value = 85
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 25
$ws.Range("F7").Value = 2.58
$ws.Range("G7").Value = 2.7
$ws.Range("H7").Value = 2.92
$ws.Range("I7").Value = 2.98
$ws.Range("J7").Value = 3.4
$ws.Range("K7").Value = 3.7
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.82
$ws.Range("W7").Value = 1.59
$ws.Range("X7").Value = 90
$ws.Range("AD7").Value = 25
$ws.Range("AL7").Value = 1000
$ws.Range("K8").Value = 5.6
$ws.Range("O8").Value = 1.21
$ws.Range("R8").Value = 1.61
$ws.Range("AB8").Value = 32
$ws.Range("AG8").Value = 29
$ws.Range("AH8").Value = 22
$ws.Range("F9").Value = 4.9
$ws.Range("H9").Value = 1.69
$ws.Range("I9").Value = 1.75
$ws.Range("J9").Value = 4.1
$ws.Range("L9").Value = 1.28
$ws.Range("Q9").Value = 1.56
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 2.34
$ws.Range("T9").Value = 1.6
$ws.Range("U9").Value = 2.54
$ws.Range("V9").Value = 2.32
$ws.Range("W9").Value = 1.22
$ws.Range("X9").Value = 55
$ws.Range("Y9").Value = 22
$ws.Range("Z9").Value = 24
$ws.Range("AB9").Value = 29
$ws.Range("AH9").Value = 24
$ws.Range("AJ9").Value = 700
$ws.Range("AL9").Value = 120
$ws.Range("G10").Value = 4.4
$ws.Range("X10").Value = 17.5
$ws.Range("AA10").Value = 20
$ws.Range("AB10").Value = 17
$ws.Range("AG10").Value = 16.5
$ws.Range("AM10").Value = 80
$ws.Range("G11").Value = 1.41
$ws.Range("L11").Value = 1.39
$ws.Range("P11").Value = 1.87
$ws.Range("Q11").Value = 1.97
$ws.Range("W11").Value = 3.4
$ws.Range("X11").Value = 14.5
$ws.Range("AC11").Value = 11.5
$ws.Range("AD11").Value = 1000
$ws.Range("AF11").Value = 7.4
$ws.Range("AH11").Value = 1000
$ws.Range("AK11").Value = 17.5
$ws.Range("F12").Value = 2.44
$ws.Range("G12").Value = 2.48
$ws.Range("T12").Value = 1.88
$ws.Range("U12").Value = 2.08
$ws.Range("W12").Value = 1.67
$ws.Range("X12").Value = 12.5
$ws.Range("Y12").Value = 12
$ws.Range("Z12").Value = 22
$ws.Range("AA12").Value = 60
$ws.Range("AE12").Value = 40
$ws.Range("AF12").Value = 14
$ws.Range("AJ12").Value = 32
$ws.Range("AK12").Value = 26
$ws.Range("AN12").Value = 24
$ws.Range("AO12").Value = 42
$ws.Range("P13").Value = 2.3
$ws.Range("Q13").Value = 1.75
$ws.Range("S13").Value = 2.9
$ws.Range("X13").Value = 22
$ws.Range("Y13").Value = 8.199999999999999
$ws.Range("AI13").Value = 44
$ws.Range("AN13").Value = 380
$ws.Range("F14").Value = 2.72
$ws.Range("G14").Value = 2.76
$ws.Range("N14").Value = 4
$ws.Range("O14").Value = 1.31
$ws.Range("P14").Value = 2.02
$ws.Range("Q14").Value = 1.96
$ws.Range("R14").Value = 1.4
$ws.Range("S14").Value = 3.45
$ws.Range("U14").Value = 2.28
$ws.Range("AA14").Value = 44
$ws.Range("AD14").Value = 12
$ws.Range("AF14").Value = 17.5
$ws.Range("H15").Value = 1.96
$ws.Range("I15").Value = 1.98
$ws.Range("O15").Value = 1.24
$ws.Range("Q15").Value = 1.72
$ws.Range("R15").Value = 1.53
$ws.Range("S15").Value = 2.8
$ws.Range("V15").Value = 2.02
$ws.Range("X15").Value = 19.5
$ws.Range("Y15").Value = 11.5
$ws.Range("AB15").Value = 18.5
$ws.Range("AF15").Value = 30
$ws.Range("AH15").Value = 15.5
$ws.Range("AL15").Value = 44
$ws.Range("F16").Value = 4.8
$ws.Range("G16").Value = 4.9
$ws.Range("O16").Value = 1.18
$ws.Range("P16").Value = 2.68
$ws.Range("Q16").Value = 1.58
$ws.Range("S16").Value = 2.4
$ws.Range("T16").Value = 1.59
$ws.Range("U16").Value = 2.62
$ws.Range("Y16").Value = 13
$ws.Range("AC16").Value = 10
$ws.Range("AH16").Value = 15.5
$ws.Range("AO16").Value = 7.2
$ws.Range("F17").Value = 2.06
$ws.Range("I17").Value = 4.4
$ws.Range("L17").Value = 1.48
$ws.Range("N17").Value = 3.05
$ws.Range("P17").Value = 1.68
$ws.Range("R17").Value = 1.23
$ws.Range("T17").Value = 1.96
$ws.Range("V17").Value = 1.29
$ws.Range("X17").Value = 12.5
$ws.Range("Y17").Value = 1000
$ws.Range("Z17").Value = 36
$ws.Range("AN17").Value = 1000
